# Add a "Save" column (H) to the s_vals sheet, mirroring the header
# formatting already used by the other header cells (column G: "sum").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (bold font, thin border, centered
# alignment) from G1 onto the new H1 header cell, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("H1").Value = "Save"

# New data row value for the Save column.
$ws.Range("H2").Value = 1
